$wb = $excel.ActiveWorkbook

# --- Sheet "NextBus3" (Worksheets index 1) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 45702.47746527778
$ws1.Cells.Item(2, 10).Value = 1
$ws1.Cells.Item(2, 15).Value = 30
$ws1.Cells.Item(3, 6).Value = 45702.47384259259
$ws1.Cells.Item(3, 12).Value = "BD"
$ws1.Cells.Item(3, 15).Value = 25
$ws1.Cells.Item(4, 6).Value = 45702.4687962963
$ws1.Cells.Item(4, 12).Value = "DD"
$ws1.Cells.Item(4, 15).Value = 18
$ws1.Cells.Item(5, 6).Value = 45702.48324074074
$ws1.Cells.Item(5, 15).Value = 39
$ws1.Cells.Item(6, 2).Value = 74
$ws1.Cells.Item(6, 3).Value = 11379
$ws1.Cells.Item(6, 4).Value = "Buona Vista Ter"
$ws1.Cells.Item(6, 5).Value = "SBST"
$ws1.Cells.Item(6, 6).Value = 45702.47726851852
$ws1.Cells.Item(6, 7).Value = 11379
$ws1.Cells.Item(6, 10).Value = 1
$ws1.Cells.Item(6, 11).Value = 64009
$ws1.Cells.Item(6, 12).Value = "DD"
$ws1.Cells.Item(6, 15).Value = 30
$ws1.Cells.Item(7, 2).Value = 61
$ws1.Cells.Item(7, 5).Value = "SMRT"
$ws1.Cells.Item(7, 6).Value = 45702.47815972222
$ws1.Cells.Item(7, 10).Value = 0
$ws1.Cells.Item(7, 11).Value = 43009
$ws1.Cells.Item(7, 12).Value = "SD"
$ws1.Cells.Item(7, 13).Value = 12109
$ws1.Cells.Item(7, 14).Value = "Opp Ngee Ann Poly"
$ws1.Cells.Item(7, 15).Value = 31
$ws1.Cells.Item(8, 3).Value = 82009
$ws1.Cells.Item(8, 4).Value = "Eunos Int"
$ws1.Cells.Item(8, 6).Value = 45702.47988425926
$ws1.Cells.Item(8, 7).Value = 82009
$ws1.Cells.Item(8, 11).Value = 22009
$ws1.Cells.Item(8, 13).Value = 12101
$ws1.Cells.Item(8, 14).Value = "Ngee Ann Poly"
$ws1.Cells.Item(8, 15).Value = 34
$ws1.Cells.Item(9, 2).Value = 154
$ws1.Cells.Item(9, 3).Value = 22009
$ws1.Cells.Item(9, 4).Value = "Boon Lay Int"
$ws1.Cells.Item(9, 5).Value = "SBST"
$ws1.Cells.Item(9, 6).Value = 45702.46974537037
$ws1.Cells.Item(9, 7).Value = 22009
$ws1.Cells.Item(9, 12).Value = "DD"
$ws1.Cells.Item(9, 13).Value = 12109
$ws1.Cells.Item(9, 14).Value = "Opp Ngee Ann Poly"
$ws1.Cells.Item(9, 15).Value = 19
$ws1.Cells.Item(10, 2).Value = 61
$ws1.Cells.Item(10, 3).Value = 43009
$ws1.Cells.Item(10, 4).Value = "Bt Batok Int"
$ws1.Cells.Item(10, 5).Value = "SMRT"
$ws1.Cells.Item(10, 6).Value = 45702.47508101852
$ws1.Cells.Item(10, 7).Value = 43009
$ws1.Cells.Item(10, 10).Value = 1
$ws1.Cells.Item(10, 11).Value = 82009
$ws1.Cells.Item(10, 15).Value = 27
$ws1.Cells.Item(11, 2).Value = 151
$ws1.Cells.Item(11, 6).Value = 45702.47321759259
$ws1.Cells.Item(11, 10).Value = 0
$ws1.Cells.Item(11, 11).Value = 16009
$ws1.Cells.Item(11, 15).Value = 24
$ws1.Cells.Item(12, 2).Value = 74
$ws1.Cells.Item(12, 3).Value = 64009
$ws1.Cells.Item(12, 4).Value = "Hougang Ctrl Int"
$ws1.Cells.Item(12, 6).Value = 45702.47288194444
$ws1.Cells.Item(12, 7).Value = 64009
$ws1.Cells.Item(12, 10).Value = 0
$ws1.Cells.Item(12, 11).Value = 11379
$ws1.Cells.Item(12, 13).Value = 12101
$ws1.Cells.Item(12, 14).Value = "Ngee Ann Poly"
$ws1.Cells.Item(12, 15).Value = 24
$ws1.Cells.Item(13, 2).Value = 52
$ws1.Cells.Item(13, 3).Value = 28009
$ws1.Cells.Item(13, 4).Value = "Jurong East Int"
$ws1.Cells.Item(13, 6).Value = 45702.47997685185
$ws1.Cells.Item(13, 7).Value = 28009
$ws1.Cells.Item(13, 9).Value = "SDA"
$ws1.Cells.Item(13, 11).Value = 53009
$ws1.Cells.Item(13, 15).Value = 34
$ws1.Cells.Item(14, 1).Value = "NextBus3"
$ws1.Cells.Item(14, 2).Value = 151
$ws1.Cells.Item(14, 3).Value = 16009
$ws1.Cells.Item(14, 4).Value = "Kent Ridge Ter"
$ws1.Cells.Item(14, 5).Value = "SBST"
$ws1.Cells.Item(14, 6).Value = 45702.47380787037
$ws1.Cells.Item(14, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(14, 7).Value = 16009
$ws1.Cells.Item(14, 8).Value = "WAB"
$ws1.Cells.Item(14, 9).Value = "SEA"
$ws1.Cells.Item(14, 10).Value = 1
$ws1.Cells.Item(14, 11).Value = 64009
$ws1.Cells.Item(14, 12).Value = "SD"
$ws1.Cells.Item(14, 13).Value = 12109
$ws1.Cells.Item(14, 14).Value = "Opp Ngee Ann Poly"
$ws1.Cells.Item(14, 15).Value = 25
$ws1.Cells.Item(15, 1).Value = "NextBus3"
$ws1.Cells.Item(15, 2).Value = 75
$ws1.Cells.Item(15, 3).Value = 10009
$ws1.Cells.Item(15, 4).Value = "Bt Merah Int"
$ws1.Cells.Item(15, 5).Value = "SMRT"
$ws1.Cells.Item(15, 6).Value = 45702.48092592593
$ws1.Cells.Item(15, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(15, 7).Value = 10009
$ws1.Cells.Item(15, 8).Value = "WAB"
$ws1.Cells.Item(15, 9).Value = "SEA"
$ws1.Cells.Item(15, 10).Value = 0
$ws1.Cells.Item(15, 11).Value = 44989
$ws1.Cells.Item(15, 12).Value = "SD"
$ws1.Cells.Item(15, 13).Value = 12109
$ws1.Cells.Item(15, 14).Value = "Opp Ngee Ann Poly"
$ws1.Cells.Item(15, 15).Value = 35

# --- Sheet "NextBus2" (Worksheets index 3) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 6).Value = 45702.47366898148
$ws3.Cells.Item(2, 12).Value = "SD"
$ws3.Cells.Item(2, 15).Value = 25
$ws3.Cells.Item(3, 6).Value = 45702.46424768519
$ws3.Cells.Item(3, 12).Value = "DD"
$ws3.Cells.Item(3, 15).Value = 11
$ws3.Cells.Item(4, 6).Value = 45702.46240740741
$ws3.Cells.Item(4, 12).Value = "DD"
$ws3.Cells.Item(4, 15).Value = 9
$ws3.Cells.Item(5, 6).Value = 45702.47521990741
$ws3.Cells.Item(5, 15).Value = 27
$ws3.Cells.Item(6, 6).Value = 45702.46640046296
$ws3.Cells.Item(6, 12).Value = "SD"
$ws3.Cells.Item(6, 15).Value = 14
$ws3.Cells.Item(7, 6).Value = 45702.46982638889
$ws3.Cells.Item(7, 10).Value = 0
$ws3.Cells.Item(7, 15).Value = 19
$ws3.Cells.Item(8, 6).Value = 45702.46943287037
$ws3.Cells.Item(8, 15).Value = 19
$ws3.Cells.Item(9, 6).Value = 45702.46556712963
$ws3.Cells.Item(9, 15).Value = 13
$ws3.Cells.Item(10, 6).Value = 45702.46450231481
$ws3.Cells.Item(10, 15).Value = 12
$ws3.Cells.Item(11, 6).Value = 45702.46461805556
$ws3.Cells.Item(11, 10).Value = 1
$ws3.Cells.Item(11, 15).Value = 12
$ws3.Cells.Item(12, 6).Value = 45702.46383101852
$ws3.Cells.Item(12, 15).Value = 11
$ws3.Cells.Item(13, 6).Value = 45702.47126157407
$ws3.Cells.Item(13, 15).Value = 21
$ws3.Cells.Item(14, 6).Value = 45702.46895833333
$ws3.Cells.Item(15, 6).Value = 45702.47390046297
$ws3.Cells.Item(15, 10).Value = 1
$ws3.Cells.Item(15, 15).Value = 25
